# Daily attendance processing - reorders the "Recorded By" (column G) list
# of names/emails for each data row: the comma-separated values are
# reversed in order (rows whose list contains "admin@admin.com", or that
# only have a single entry, are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }

    $s = [string]$raw
    if ($s -eq "") { continue }
    if ($s.Contains("admin@admin.com")) { continue }
    if (-not $s.Contains(",")) { continue }

    $parts = $s.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $n = $trimmed.Length
    if ($n -le 1) { continue }

    $rev = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $rev += $trimmed[$i]
    }

    $joined = [string]::Join(", ", $rev)
    $cell.Value = $joined
}
